$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 41 (ALC)
$ws.Range("H41").Value = 5007.5
$ws.Range("I41").Value = 267
$ws.Range("J41").Value = 8289.385
$ws.Range("K41").Value = 267
$ws.Range("L41").Value = 8289.385
$ws.Range("M41").Value = 173
$ws.Range("N41").Value = -9169.385

# row 76 (ALC)
$ws.Range("H76").Value = 8991.041999999999
$ws.Range("I76").Value = 13868.1
$ws.Range("K76").Value = 13868.1
$ws.Range("M76").Value = -13553.1

# row 79 (ALC)
$ws.Range("H79").Value = 8991.041999999999
$ws.Range("I79").Value = 13868.1
$ws.Range("K79").Value = 13868.1
$ws.Range("M79").Value = -12776.1

# row 107 (ALC)
$ws.Range("H107").Value = 748.9
$ws.Range("I107").Value = 336.125
$ws.Range("J107").Value = 2400
$ws.Range("K107").Value = 336.125
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = 1583.875
$ws.Range("N107").Value = -6240

# row 132 (ALC)
$ws.Range("H132").Value = 2204.6323
$ws.Range("I132").Value = 1201.1818
$ws.Range("K132").Value = 3603.5454
$ws.Range("M132").Value = -1073.5454

# row 136 (ALC)
$ws.Range("H136").Value = 33569.6
$ws.Range("J136").Value = 33569.6
$ws.Range("L136").Value = 33569.6
$ws.Range("N136").Value = -43769.6

# row 138 (ALC)
$ws.Range("H138").Value = 1772.0847
$ws.Range("I138").Value = 1284.5862
$ws.Range("J138").Value = 2243.3333
$ws.Range("K138").Value = 3853.7586
$ws.Range("L138").Value = 6729.999899999999
$ws.Range("M138").Value = 1286.2414
$ws.Range("N138").Value = -17009.9999

# row 141 (ALC)
$ws.Range("H141").Value = 1509.8889
$ws.Range("I141").Value = 1448.625
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 4345.875
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 834.125
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
# row 55 (ARM)
$ws.Range("H55").Value = 23500

# row 102 (ARM)
$ws.Range("H102").Value = 3052.975
$ws.Range("I102").Value = 1552.8387
$ws.Range("K102").Value = 1552.8387
$ws.Range("M102").Value = 69.16129999999998

# row 110 (ARM)
$ws.Range("H110").Value = 1253.6897
$ws.Range("I110").Value = 1174.7727
$ws.Range("K110").Value = 1174.7727
$ws.Range("M110").Value = 870.2273

# row 122 (ARM)
$ws.Range("H122").Value = 2675.7778
$ws.Range("I122").Value = 2699.7
$ws.Range("J122").Value = 2645.875
$ws.Range("K122").Value = 8099.099999999999
$ws.Range("L122").Value = 7937.625
$ws.Range("M122").Value = -5649.099999999999
$ws.Range("N122").Value = -12837.625

$ws = $wb.Worksheets.Item("BSM")
# row 20 (BSM)
$ws.Range("H20").Value = 1749.7037
$ws.Range("I20").Value = 1757.4286
$ws.Range("J20").Value = 1741.3846
$ws.Range("K20").Value = 1757.4286
$ws.Range("L20").Value = 1741.3846
$ws.Range("M20").Value = -1510.4286
$ws.Range("N20").Value = -2235.3846

# row 86 (BSM)
$ws.Range("H86").Value = 6023.696
$ws.Range("I86").Value = 7102
$ws.Range("J86").Value = 5330.5
$ws.Range("K86").Value = 7102
$ws.Range("L86").Value = 5330.5
$ws.Range("M86").Value = -5979
$ws.Range("N86").Value = -7576.5

# row 87 (BSM)
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").ClearContents()

# row 89 (BSM)
$ws.Range("H89").Value = 6023.696
$ws.Range("I89").Value = 7102
$ws.Range("J89").Value = 5330.5
$ws.Range("K89").Value = 35510
$ws.Range("L89").Value = 26652.5
$ws.Range("M89").Value = -29894
$ws.Range("N89").Value = -37884.5

# row 90 (BSM)
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").ClearContents()

# row 99 (BSM)
$ws.Range("H99").Value = 3726.5557
$ws.Range("I99").Value = 5100.6665
$ws.Range("J99").Value = 978.3333
$ws.Range("K99").Value = 5100.6665
$ws.Range("L99").Value = 978.3333
$ws.Range("M99").Value = -3602.6665
$ws.Range("N99").Value = -3974.3333

$ws = $wb.Worksheets.Item("CRP")
# row 7 (CRP)
$ws.Range("H7").Value = 53.92857
$ws.Range("I7").Value = 36
$ws.Range("J7").Value = 86.2
$ws.Range("K7").Value = 36
$ws.Range("L7").Value = 86.2
$ws.Range("M7").Value = 77
$ws.Range("N7").Value = -312.2

# row 14 (CRP)
$ws.Range("H14").Value = 950
$ws.Range("J14").Value = 950
$ws.Range("L14").Value = 950
$ws.Range("N14").Value = -1290

# row 107 (CRP)
$ws.Range("H107").Value = 332.85715
$ws.Range("I107").Value = 262.33334
$ws.Range("J107").Value = 459.8
$ws.Range("K107").Value = 262.33334
$ws.Range("L107").Value = 459.8
$ws.Range("M107").Value = 1657.66666
$ws.Range("N107").Value = -4299.8

$ws = $wb.Worksheets.Item("GSM")
# row 70 (GSM)
$ws.Range("H70").Value = 5530.472
$ws.Range("I70").Value = 3975.4285
$ws.Range("J70").Value = 8843.392
$ws.Range("K70").Value = 3975.4285
$ws.Range("L70").Value = 8843.392
$ws.Range("M70").Value = -3705.4285
$ws.Range("N70").Value = -9383.392

# row 73 (GSM)
$ws.Range("H73").Value = 5530.472
$ws.Range("I73").Value = 3975.4285
$ws.Range("J73").Value = 8843.392
$ws.Range("K73").Value = 3975.4285
$ws.Range("L73").Value = 8843.392
$ws.Range("M73").Value = -3039.4285
$ws.Range("N73").Value = -10715.392

# row 80 (GSM)
$ws.Range("H80").Value = 4458.3335
$ws.Range("I80").Value = 5431.304
$ws.Range("J80").Value = 2736.923
$ws.Range("K80").Value = 5431.304
$ws.Range("L80").Value = 2736.923
$ws.Range("M80").Value = -4433.304
$ws.Range("N80").Value = -4732.923

# row 83 (GSM)
$ws.Range("H83").Value = 4458.3335
$ws.Range("I83").Value = 5431.304
$ws.Range("J83").Value = 2736.923
$ws.Range("K83").Value = 27156.52
$ws.Range("L83").Value = 13684.615
$ws.Range("M83").Value = -22164.52
$ws.Range("N83").Value = -23668.615

# row 122 (GSM)
$ws.Range("H122").Value = 1160
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -10300

# row 141 (GSM)
$ws.Range("H141").Value = 36235.75
$ws.Range("J141").Value = 43462.668
$ws.Range("L141").Value = 43462.668
$ws.Range("N141").Value = -53822.668

$ws = $wb.Worksheets.Item("LTW")
# row 55 (LTW)
$ws.Range("H55").Value = 715.2857
$ws.Range("I55").Value = 426.5
$ws.Range("J55").Value = 1100.3334
$ws.Range("K55").Value = 426.5
$ws.Range("L55").Value = 1100.3334
$ws.Range("M55").Value = -253.5
$ws.Range("N55").Value = -1446.3334

$ws = $wb.Worksheets.Item("WVR")
# row 54 (WVR)
$ws.Range("H54").Value = 14800
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 14800
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 14800
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -15840

# row 140 (WVR)
$ws.Range("H140").Value = 26754.857
$ws.Range("J140").Value = 26754.857
$ws.Range("L140").Value = 26754.857
$ws.Range("N140").Value = -37114.857

# row 141 (WVR)
$ws.Range("H141").Value = 30143.334
$ws.Range("J141").Value = 30143.334
$ws.Range("L141").Value = 30143.334
$ws.Range("N141").Value = -40503.334
